# #3456 updated PM Property ID
#
# The "Portfolio Manager Building ID" column (column B) on the "BPS Data"
# sheet gets new PM Property IDs for rows 2-10, and the sheet's active
# selection moves from A2:A10 to B16:C24 (reflecting where the user was
# last working in the real edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Portfolio Manager Building ID (column B) values, rows 2-10.
$ws.Range("B2").Value  = 22178843
$ws.Range("B3").Value  = 22178844
$ws.Range("B4").Value  = 22178845
$ws.Range("B5").Value  = 22178846
$ws.Range("B6").Value  = 22178847
$ws.Range("B7").Value  = 22178848
$ws.Range("B8").Value  = 22178849
$ws.Range("B9").Value  = 22178850
$ws.Range("B10").Value = 22178851

# Move the sheet's active selection to match the saved view state.
$ws.Range("B16:C24").Select()
